$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the
# "Technologies, Tools, and Resources Used" heading paragraph. It needs to
# move to the end of the new last bullet under "Tasks Undertaken". Remove it
# from its current spot first; we'll re-add it in the right place below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# Locate the "Copied project from Task 14..." bullet (currently the last
# bullet under "Tasks Undertaken") by scanning paragraphs, so we don't
# depend on hard-coded paragraph indices.
$copiedPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Copied project from Task 14*") {
        $copiedPara = $p
        break
    }
}

# Add the new bullet: "Set up a pre-set patrol path ..."
$copiedPara.Range.InsertParagraphAfter() | Out-Null
$newPara1 = $copiedPara.Next()
$newPara1.Range.Text = "Set up a pre-set patrol path and changed the soldier and target" + [char]0x2019 + "s starting positions to be the first point in the path and the centre of the simulation space respectively."

# Add the new bullet: "Altered soldier agent to expand its field of view ..."
$newPara1.Range.InsertParagraphAfter() | Out-Null
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "Altered soldier agent to expand its field of view and remove gaps in  said field of view."

# Re-create the "_GoBack" bookmark at the end of the new last bullet's text
# (immediately after the text, before the paragraph mark), matching its
# original position relative to the text it used to follow.
$endRange = $newPara2.Range.Duplicate
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
